$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 82; existing rows 82..196 shift down to 83..197.
$ws.Rows.Item(82).Insert()

$r = 82
$ws.Cells.Item($r,1).Value = 5
$ws.Cells.Item($r,2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item($r,3).Value = "Maule"
$ws.Cells.Item($r,4).Value = 44482
$ws.Cells.Item($r,5).Value = 7
$ws.Cells.Item($r,6).Value = 100114013
$ws.Cells.Item($r,7).Value = "Zanahoria"
$ws.Cells.Item($r,8).Value = "Sin especificar"
$ws.Cells.Item($r,9).Value = "Primera"
$ws.Cells.Item($r,10).Value = 500
$ws.Cells.Item($r,11).Value = 8000
$ws.Cells.Item($r,12).Value = 8000
$ws.Cells.Item($r,13).Value = 8000
$ws.Cells.Item($r,14).Value = "$/saco 20 kilos"
$ws.Cells.Item($r,15).Value = "Región de Ñuble"
$ws.Cells.Item($r,16).Value = 400
$ws.Cells.Item($r,17).Value = 20
$ws.Cells.Item($r,18).Value = "Hortaliza"
